# Updates the cryptocurrency price/volume table (Sheet1, columns B:E) to
# reflect the latest scrape: refreshed Price (D) and Volume/1h (E) figures
# for every coin, plus two pairs of rows (29/30 and 41/42) that swapped
# ranking order, bringing a different coin's Name/Link/Price/Volume into
# each of those rows.
#
# Price-column values often look numeric ("234.05", "0.0760", ...) but the
# source data stores them as plain text (so things like trailing zeros -
# "15.00", "0.0760" - survive). A bare Range.Value assignment would let
# Excel's automatic type-detection coerce those into real numbers and
# silently drop the trailing zeros, so every Price (column D) write below
# is apostrophe-prefixed to force text entry, exactly like a user typing
# '234.05 into a cell. The leading apostrophe is not part of the stored
# value - Excel strips it and simply marks the cell as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'37.283.90"
$ws.Range("E2").Value = "  +0.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.062.22"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'234.05"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  +1.76%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'56.61"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.63%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0760"
$ws.Range("E10").Value = "  +0.19%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.56%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "'2.366.65"
$ws.Range("E12").Value = "  -0.10%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "'14.59"
$ws.Range("E13").Value = "  +1.03%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  -2.52%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.36%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "'5.10"
$ws.Range("E16").Value = "  -2.41%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'2.062.16"
$ws.Range("E17").Value = "  -0.12%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'37.239.83"
$ws.Range("E18").Value = "  +0.09%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +6.85%  "

# Row 20 - Litecoin
$ws.Range("D20").Value = "'69.42"
$ws.Range("E20").Value = "  +1.69%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "'0.0₃0808"
$ws.Range("E21").Value = "  +0.05%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'225.99"
$ws.Range("E22").Value = "  +1.15%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.04%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +1.40%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -1.43%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'166.04"
$ws.Range("E26").Value = "  +2.07%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +4.79%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "'8.75"
$ws.Range("E28").Value = "  -0.84%  "

# Row 29 - now EthereumClassic (swapped with row 30's former occupant)
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'18.99"
$ws.Range("E29").Value = "  -1.08%  "

# Row 30 - now Kaspa (swapped with row 29's former occupant)
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.126"
$ws.Range("E30").Value = "  -2.74%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.117"
$ws.Range("E31").Value = "  -0.67%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.25%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -0.82%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +3.74%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -0.86%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  +0.05%  "

# Row 37 - WEMIXToken
$ws.Range("D37").Value = "'1.76"
$ws.Range("E37").Value = "  -0.87%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -2.78%  "

# Row 39 - THORChain
$ws.Range("D39").Value = "'5.67"
$ws.Range("E39").Value = "  -4.48%  "

# Row 40 - HuobiToken
$ws.Range("E40").Value = "  -0.37%  "

# Row 41 - now Maker (swapped with row 42's former occupant)
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.461.73"
$ws.Range("E41").Value = "  -0.57%  "

# Row 42 - now Aave (swapped with row 41's former occupant)
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'95.97"
$ws.Range("E42").Value = "  +1.83%  "

# Row 43 - Cronos
$ws.Range("D43").Value = "'0.0935"
$ws.Range("E43").Value = "  -2.64%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +1.90%  "

# Row 45 - TrustWalletToken
$ws.Range("E45").Value = "  +2.96%  "

# Row 46 - FTXToken
$ws.Range("D46").Value = "'4.20"
$ws.Range("E46").Value = "  -5.57%  "

# Row 47 - ARBITRUM
$ws.Range("E47").Value = "  -0.17%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "'15.00"
$ws.Range("E48").Value = "  -6.44%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  -0.08%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +1.10%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "'2.252.89"
$ws.Range("E51").Value = "  -0.14%  "
